# Updated cryptos list on Sat Aug 24 03:25:33 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# every coin row, and fixes two rows where the scraped ranking had swapped
# two coins' positions (Dai/Polygon at rows 23-24, Aave/Mantle at rows
# 46-47) - each row keeps the same row number but gets the other coin's
# name/link/price/volume.
#
# Column D ("Price") cells are stored as literal text in the sheet (e.g.
# "0.390", "63.548.71" - note some use '.' as a thousands separator, which
# isn't valid numeric syntax anyway), so each Price write temporarily forces
# the cell to Text format, assigns the literal string, then restores the
# cell's default ("Normal") style - this prevents Excel's normal type
# inference from silently reinterpreting the text as a number (which would
# both change its stored type and, for values like "0.390", drop the
# significant trailing zero).
#
# Column E ("Volume(1h)") values keep their padding spaces (e.g.
# "  +4.95%  ") and column B/C values are plain names/URLs, so neither
# needs that treatment - a direct .Value assignment round-trips them as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.548.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.724.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.74%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.607'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.752.70'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.72'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.11%  '
$ws.Range('E11').Value = '  +6.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.390'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.160'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.213.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.29'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.746.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.30%  '
$ws.Range('E17').Value = '  +6.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.754.24'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.97'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '360.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.536'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.995'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.35%  '
$ws.Range('E26').Value = '  +4.63%  '
$ws.Range('E27').Value = '  +5.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  +11.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '173.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.52'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('E36').Value = '  +7.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.60%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.81'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.997'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +15.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '344.91'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.64%  '
$ws.Range('E41').Value = '  +5.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.09'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.54'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.69'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.83'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.94%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.649'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.15%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '139.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.30%  '
$ws.Range('E48').Value = '  +5.31%  '
$ws.Range('E49').Value = '  +4.65%  '
$ws.Range('E50').Value = '  +0.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.996'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.32%  '
